$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell while forcing the
# Text number format so Excel does not silently reinterpret numeric-
# looking strings (e.g. "0.998", "1.00") as actual numbers. The cell
# style is then reset back to "Normal" so no stray formatting is left
# behind relative to the original (unstyled) cells.
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '66.505.48'
Set-TextValue $ws.Range('E2') '  -2.71%  '
Set-TextValue $ws.Range('D3') '3.564.51'
Set-TextValue $ws.Range('E3') '  -3.10%  '
Set-TextValue $ws.Range('E4') '  -0.39%  '
Set-TextValue $ws.Range('D5') '586.65'
Set-TextValue $ws.Range('E5') '  -1.07%  '
Set-TextValue $ws.Range('D6') '181.35'
Set-TextValue $ws.Range('E6') '  +0.65%  '
Set-TextValue $ws.Range('D7') '3.560.16'
Set-TextValue $ws.Range('E7') '  -3.08%  '
Set-TextValue $ws.Range('D8') '0.608'
Set-TextValue $ws.Range('E8') '  -3.40%  '
Set-TextValue $ws.Range('D9') '0.998'
Set-TextValue $ws.Range('E9') '  -0.19%  '
Set-TextValue $ws.Range('D10') '0.668'
Set-TextValue $ws.Range('E10') '  -6.46%  '
Set-TextValue $ws.Range('B11') 'Avalanche'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D11') '53.62'
Set-TextValue $ws.Range('E11') '  -3.80%  '
Set-TextValue $ws.Range('B12') 'Dogecoin'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D12') '0.143'
Set-TextValue $ws.Range('E12') '  -10.63%  '
Set-TextValue $ws.Range('D13') '0.0000251'
Set-TextValue $ws.Range('E13') '  -14.09%  '
Set-TextValue $ws.Range('D14') '9.74'
Set-TextValue $ws.Range('E14') '  -8.40%  '
Set-TextValue $ws.Range('D15') '4.127.56'
Set-TextValue $ws.Range('E15') '  -3.35%  '
Set-TextValue $ws.Range('D16') '3.558.43'
Set-TextValue $ws.Range('E16') '  -3.38%  '
Set-TextValue $ws.Range('E17') '  -0.49%  '
Set-TextValue $ws.Range('D18') '18.28'
Set-TextValue $ws.Range('E18') '  -5.52%  '
Set-TextValue $ws.Range('D19') '66.221.91'
Set-TextValue $ws.Range('E19') '  -2.96%  '
Set-TextValue $ws.Range('D20') '12.07'
Set-TextValue $ws.Range('E20') '  -5.84%  '
Set-TextValue $ws.Range('D21') '1.05'
Set-TextValue $ws.Range('E21') '  -6.63%  '
Set-TextValue $ws.Range('D22') '391.84'
Set-TextValue $ws.Range('E22') '  -4.53%  '
Set-TextValue $ws.Range('D23') '4.29'
Set-TextValue $ws.Range('E23') '  -6.55%  '
Set-TextValue $ws.Range('D24') '84.32'
Set-TextValue $ws.Range('E24') '  -4.82%  '
Set-TextValue $ws.Range('D25') '2.86'
Set-TextValue $ws.Range('E25') '  -5.18%  '
Set-TextValue $ws.Range('D26') '12.25'
Set-TextValue $ws.Range('E26') '  -3.42%  '
Set-TextValue $ws.Range('D27') '6.03'
Set-TextValue $ws.Range('E27') '  -0.62%  '
Set-TextValue $ws.Range('D28') '10.20'
Set-TextValue $ws.Range('E28') '  -5.35%  '
Set-TextValue $ws.Range('D29') '3.56'
Set-TextValue $ws.Range('E29') '  -8.30%  '
Set-TextValue $ws.Range('D30') '8.88'
Set-TextValue $ws.Range('E30') '  -7.45%  '
Set-TextValue $ws.Range('D31') '30.96'
Set-TextValue $ws.Range('E31') '  -5.67%  '
Set-TextValue $ws.Range('D32') '6.71'
Set-TextValue $ws.Range('E32') '  -8.78%  '
Set-TextValue $ws.Range('B33') 'Cosmos'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D33') '11.89'
Set-TextValue $ws.Range('E33') '  -4.26%  '
Set-TextValue $ws.Range('B34') 'OKB'
Set-TextValue $ws.Range('C34') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D34') '65.10'
Set-TextValue $ws.Range('E34') '  +0.64%  '
Set-TextValue $ws.Range('D35') '603.30'
Set-TextValue $ws.Range('E35') '  -0.23%  '
Set-TextValue $ws.Range('E36') '  -6.76%  '
Set-TextValue $ws.Range('D37') '41.25'
Set-TextValue $ws.Range('E37') '  -4.64%  '
Set-TextValue $ws.Range('E38') '  +0.07%  '
Set-TextValue $ws.Range('E39') '  -0.32%  '
Set-TextValue $ws.Range('D40') '0.371'
Set-TextValue $ws.Range('E40') '  -7.75%  '
Set-TextValue $ws.Range('D41') '0.0₃0737'
Set-TextValue $ws.Range('E41') '  -15.75%  '
Set-TextValue $ws.Range('D42') '0.129'
Set-TextValue $ws.Range('E42') '  -6.59%  '
Set-TextValue $ws.Range('D43') '2.900.91'
Set-TextValue $ws.Range('E43') '  +6.50%  '
Set-TextValue $ws.Range('D44') '2.76'
Set-TextValue $ws.Range('E44') '  -8.85%  '
Set-TextValue $ws.Range('D45') '0.0405'
Set-TextValue $ws.Range('E45') '  -8.17%  '
Set-TextValue $ws.Range('D46') '2.40'
Set-TextValue $ws.Range('E46') '  -9.25%  '
Set-TextValue $ws.Range('D47') '0.130'
Set-TextValue $ws.Range('E47') '  -3.97%  '
Set-TextValue $ws.Range('D48') '3.03'
Set-TextValue $ws.Range('E48') '  -2.23%  '
Set-TextValue $ws.Range('B49') 'Monero'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D49') '135.67'
Set-TextValue $ws.Range('E49') '  -3.14%  '
Set-TextValue $ws.Range('B50') 'WEMIXToken'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D50') '2.49'
Set-TextValue $ws.Range('E50') '  -8.68%  '
Set-TextValue $ws.Range('D51') '8.20'
Set-TextValue $ws.Range('E51') '  -8.92%  '
